$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 13 ("Total Volunteer Hours"),
# shifting rows 13-19 down to 14-20, to hold the new "Volunteers" entry.
$ws.Rows.Item(13).Insert()

# The inserted row has no formatting yet - copy the look of the row
# directly above it (the "Number of Volunteers" row) so the new row
# matches the existing label/value styling exactly.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# Fill in the new row's content.
$ws.Range("A13").Value = "Volunteers"
$ws.Range("B13").Value = "Amrit Manhas, Greg Pikatis, Alobo Dreok"

# Widen column B to fit the new, longer volunteer-name text.
$ws.Columns.Item(2).ColumnWidth = 36.4375

# Match the saved selection from the edited workbook.
$ws.Range("D16").Select()
